# Update sock PIC default product weight: 0.4 -> 0.04
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PIC default product weights")

$ws.Range("C7").Value = 0.04

# Set column widths for columns A and B (as adjusted by the author)
$ws.Columns.Item(1).ColumnWidth = 28.333333333333336
$ws.Columns.Item(2).ColumnWidth = 25.833333333333336

# Move the active selection to C8
[void]$ws.Range("C8").Select()

# Reposition the workbook window (cosmetic window move recorded by Excel)
$win = $wb.Windows.Item(1)
$win.Left = 5200
$win.Top = 3100
